$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data source row (row 3): "silverarrow" pre-production environment
$ws.Range("B3").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("A3").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("C3").Value = "su"
$ws.Range("D3").Value = "silverarrow"
$ws.Range("E3").Value = "'04104013020"
$ws.Range("F3").Value = "ABM Cláusula Ajuste"
$ws.Range("G3").Value = 10

$ws.Hyperlinks.Add($ws.Range("B3"), "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do")
$ws.Range("B3").Style = "Hipervínculo"

# Leave the selection on the newly-entered row, like after tabbing through it
$ws.Range("F3:G3").Select() | Out-Null
